# Rename the three inline logo pictures in the document's headers/footers.
#
#   * Pearson logo in the "default" footer  (footer2.xml, docPr id="2")
#         image1.png -> image2.png
#   * Pearson logo in the "first page" footer (footer1.xml, docPr id="3")
#         image1.png -> image2.png
#   * BTec logo in the "first page" header  (header1.xml, docPr id="1")
#         image2.jpg -> image1.jpg
#
# The pictures are inline drawings living in the headers/footers, so they
# are reached through Section.Headers / Section.Footers (wdHeaderFooterPrimary
# == 1, wdHeaderFooterFirstPage == 2) rather than Document.InlineShapes
# (which only covers the main story).

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footers: Pearson logo -------------------------------------------------
$defaultFooterShape = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$defaultFooterShape.Name = "image2.png"

$firstFooterShape = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$firstFooterShape.Name = "image2.png"

# --- Header: BTec logo -------------------------------------------------
$firstHeaderShape = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$firstHeaderShape.Name = "image1.jpg"
